$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 - this shifts the existing rows 36..43
# down to 37..44, preserving all of their data/formatting.
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new record.
$ws.Cells.Item(36, 1).Value = 8
$ws.Cells.Item(36, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = 45204
$ws.Cells.Item(36, 5).Value = 4
$ws.Cells.Item(36, 6).Value = 100112013
$ws.Cells.Item(36, 7).Value = "Alcachofa"
$ws.Cells.Item(36, 8).Value = "Española"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 500
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 9000
$ws.Cells.Item(36, 13).Value = 8500
$ws.Cells.Item(36, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(36, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value = 283
$ws.Cells.Item(36, 17).Value = 30
$ws.Cells.Item(36, 18).Value = "Hortaliza"
